# Added Enquiry to remaining Params
# Populate the three new "Enquiry" description cells on Sheet1 (Z2, AA2, Z3)
# with matching styling (thin border, 7.5pt Calibri, wrapped + vertically
# centered text) and size the two new columns to fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New cell content -------------------------------------------------
$ws.Range("Z2").Value  = "Account Number of the Bank for payments effected by the Client through Bank"
$ws.Range("AA2").Value = "Banking mode of payment of premiums"
$ws.Range("Z3").Value  = "Group Ledger Code for the payments effected by the Insurer"

# --- Formatting for the new cells --------------------------------------
$rng = $ws.Range("Z2:AA3")
$rng.Font.Size = 7.5
$rng.VerticalAlignment = -4108   # xlCenter
$rng.WrapText = $true
$rng.Borders.LineStyle = 1       # xlContinuous
$rng.Borders.Weight = 2          # xlThin

# --- Column widths for the new columns ---------------------------------
$ws.Columns.Item("Z").ColumnWidth = 54.6
$ws.Columns.Item("AA").ColumnWidth = 55.92

# --- Selection / view state ---------------------------------------------
$ws.Range("Z2:AA3").Select()
$excel.ActiveWindow.ScrollColumn = 17
